# #5: property aircraft done
# The "property_category" column on the 建物 (building) sheet was
# incorrectly recorded as "land" - fix it to "building".
# The "property_category" column on the 汽車 (car/vehicle) sheet was
# also incorrectly recorded as "land" - fix it to "car".

$wb = $excel.ActiveWorkbook

$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2").Value = "building"

$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
